# The source data for sheet "foo" had an extra color value ("orange") in
# cell D6 that is no longer present in the edited workbook. Clearing this
# cell also drops "orange" from the shared-strings table once the file is
# saved, which is what the target diff shows (uniqueCount 13 -> 12 and all
# later shared-string indices shifting down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("foo")

$ws.Range("D6").ClearContents()
